$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 148 (shifts existing rows 148-214 down to 149-215).
$ws.Rows("148:148").Insert()

# The row that was previously row 148 is now row 149; duplicate its contents
# into the newly-inserted (blank) row 148 so we start from an identical record.
$ws.Range("A149:R149").Copy($ws.Range("A148:R148"))

# Now overwrite the specific fields that differ for the new record.
$ws.Range("D148").Value2 = 44813
$ws.Range("K148").Value2 = 3000
$ws.Range("L148").Value2 = 3000
$ws.Range("M148").Value2 = 3000
$ws.Range("P148").Value2 = 750
